$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 4: drop the "NO HP CS" and "KATEGORI" columns ---
# Before: B4=NO, C4=KURIR PICK UP, D4=KODE RESI, E4=NAMA CS, F4=NO HP CS,
#         G4=NO HP SELLER, H4=HARGA, I4=ONGKIR, J4=KATEGORI, K4=KETERANGAN
# After:  B4=NO, C4=KURIR PICK UP, D4=KODE RESI, E4=NAMA CS, F4=NO HP SELLER,
#         G4=HARGA, H4=ONGKIR, I4=KETERANGAN
$ws.Range("F4").Value = "NO HP SELLER"
$ws.Range("G4").Value = "HARGA"
$ws.Range("H4").Value = "ONGKIR"
$ws.Range("I4").Value = "KETERANGAN"
$ws.Range("J4").Clear()
$ws.Range("K4").Clear()

# --- Un-merge the title/banner ranges before the column layout changes ---
$ws.Range("K2").UnMerge()
$ws.Range("K3").UnMerge()
$ws.Range("I8").UnMerge()

# --- Column width of the (now) HARGA column ---
$ws.Range("G1").ColumnWidth = 11.6

# --- Drop the now-unused column width definitions for the removed columns ---
$ws.Range("J1:K1").EntireColumn.Delete()

# --- Re-merge the title/banner ranges to match the narrower table ---
$ws.Range("B2:J2").Merge()
$ws.Range("B3:J3").Merge()
$ws.Range("B8:H8").Merge()
